$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

$ws.Range('E2').Value = '2026-02-14 06:18:37'
$ws.Range('G2').Value = '121 cm'
$ws.Range('I2').Value = '9.3 mm'
$ws.Range('N2').Value = '-1.1 °C 5:42 TU'
$ws.Range('E3').Value = '2026-02-14 06:18:40'
$ws.Range('I3').Value = '4.3 mm'
$ws.Range('N3').Value = '-5.3 °C 5:53 TU'
$ws.Range('E4').Value = '2026-02-14 06:18:42'
$ws.Range('H4').NumberFormat = "@"
$ws.Range('H4').Value = '83%'
$ws.Range('J4').Value = '989.8 hPa'
$ws.Range('N4').Value = '6.4 °C 5:57 TU'
$ws.Range('O4').Value = '8.4 °C'
$ws.Range('E5').Value = '2026-02-14 06:18:45'
$ws.Range('I5').Value = '7.3 mm'
$ws.Range('N5').Value = '-5.1 °C 5:58 TU'
$ws.Range('E6').Value = '2026-02-14 06:18:47'
$ws.Range('J6').Value = '989.8 hPa'
$ws.Range('N6').Value = '6.3 °C 5:55 TU'
$ws.Range('O6').Value = '7.2 °C'
$ws.Range('E7').Value = '2026-02-14 06:18:49'
$ws.Range('J7').Value = '990.1 hPa'
$ws.Range('E8').Value = '2026-02-14 06:18:52'
$ws.Range('J8').Value = '989.7 hPa'
$ws.Range('N8').Value = '6.5 °C 5:30 TU'
$ws.Range('E9').Value = '2026-02-14 06:18:55'
$ws.Range('H9').NumberFormat = "@"
$ws.Range('H9').Value = '60%'
$ws.Range('O9').Value = '11.4 °C'
$ws.Range('E10').Value = '2026-02-14 06:18:57'
$ws.Range('N10').Value = '6.1 °C 5:41 TU'
$ws.Range('E11').Value = '2026-02-14 06:18:59'
$ws.Range('H11').NumberFormat = "@"
$ws.Range('H11').Value = '92%'
$ws.Range('M11').Value = '7.4 °C 5:46 TU'
$ws.Range('O11').Value = '3.3 °C'
$ws.Range('E12').Value = '2026-02-14 06:19:02'
$ws.Range('H12').NumberFormat = "@"
$ws.Range('H12').Value = '61%'
$ws.Range('O12').Value = '11.9 °C'
$ws.Range('E13').Value = '2026-02-14 06:19:04'
$ws.Range('H13').NumberFormat = "@"
$ws.Range('H13').Value = '87%'
$ws.Range('J13').Value = '991.4 hPa'
$ws.Range('O13').Value = '2.1 °C'
$ws.Range('E14').Value = '2026-02-14 06:19:07'
$ws.Range('H14').NumberFormat = "@"
$ws.Range('H14').Value = '65%'
$ws.Range('E15').Value = '2026-02-14 06:19:09'
$ws.Range('H15').NumberFormat = "@"
$ws.Range('H15').Value = '63%'
$ws.Range('E16').Value = '2026-02-14 06:19:12'
$ws.Range('E17').Value = '2026-02-14 06:19:14'
$ws.Range('O17').Value = '0.8 °C'
$ws.Range('E18').Value = '2026-02-14 06:19:17'
$ws.Range('J18').Value = '990.0 hPa'
$ws.Range('N18').Value = '6.4 °C 5:42 TU'
$ws.Range('O18').Value = '7.4 °C'
$ws.Range('E19').Value = '2026-02-14 06:19:19'
$ws.Range('N19').Value = '2.6 °C 5:46 TU'
$ws.Range('E20').Value = '2026-02-14 06:19:22'
$ws.Range('G20').Value = '119 cm'
$ws.Range('I20').Value = '1.3 mm'
$ws.Range('L20').Value = '79.2 km/h - 335º 5:30 TU'
$ws.Range('N20').Value = '-5.4 °C 5:58 TU'
$ws.Range('E21').Value = '2026-02-14 06:19:24'
$ws.Range('H21').NumberFormat = "@"
$ws.Range('H21').Value = '95%'
$ws.Range('J21').Value = '992.5 hPa'
$ws.Range('O21').Value = '1.2 °C'
$ws.Range('E22').Value = '2026-02-14 06:19:27'
$ws.Range('I22').Value = '0.3 mm'
$ws.Range('E23').Value = '2026-02-14 06:19:29'
$ws.Range('H23').NumberFormat = "@"
$ws.Range('H23').Value = '86%'
$ws.Range('I23').Value = '10.6 mm'
$ws.Range('N23').Value = '-5.9 °C 5:37 TU'
$ws.Range('E24').Value = '2026-02-14 06:19:32'
$ws.Range('J24').Value = '993.8 hPa'
$ws.Range('O24').Value = '7.2 °C'
$ws.Range('E25').Value = '2026-02-14 06:19:34'
$ws.Range('I25').Value = '19.6 mm'
$ws.Range('E26').Value = '2026-02-14 06:19:37'
$ws.Range('E27').Value = '2026-02-14 06:19:39'
$ws.Range('H27').NumberFormat = "@"
$ws.Range('H27').Value = '83%'
$ws.Range('N27').Value = '-3.5 °C 5:59 TU'
$ws.Range('E28').Value = '2026-02-14 06:19:42'
$ws.Range('J28').Value = '990.3 hPa'
$ws.Range('N28').Value = '4.5 °C 5:59 TU'
$ws.Range('E29').Value = '2026-02-14 06:19:44'
$ws.Range('E30').Value = '2026-02-14 06:19:47'
$ws.Range('H30').NumberFormat = "@"
$ws.Range('H30').Value = '73%'
$ws.Range('J30').Value = '989.4 hPa'
$ws.Range('O30').Value = '10.1 °C'
$ws.Range('E31').Value = '2026-02-14 06:19:49'
$ws.Range('H31').NumberFormat = "@"
$ws.Range('H31').Value = '78%'
$ws.Range('J31').Value = '988.9 hPa'
$ws.Range('L31').Value = '121.7 km/h - 331º 5:40 TU'
$ws.Range('E32').Value = '2026-02-14 06:19:52'
$ws.Range('H32').NumberFormat = "@"
$ws.Range('H32').Value = '97%'
$ws.Range('N32').Value = '2.0 °C 5:38 TU'
$ws.Range('O32').Value = '2.5 °C'
$ws.Range('E33').Value = '2026-02-14 06:19:54'
$ws.Range('J33').Value = '990.1 hPa'
$ws.Range('E34').Value = '2026-02-14 06:19:57'
$ws.Range('H34').NumberFormat = "@"
$ws.Range('H34').Value = '75%'
$ws.Range('N34').Value = '-2.2 °C 5:53 TU'
$ws.Range('O34').Value = '-1.4 °C'
$ws.Range('E35').Value = '2026-02-14 06:19:59'
$ws.Range('H35').NumberFormat = "@"
$ws.Range('H35').Value = '80%'
$ws.Range('J35').Value = '995.5 hPa'
$ws.Range('N35').Value = '1.8 °C 5:50 TU'
$ws.Range('O35').Value = '2.6 °C'
$ws.Range('E36').Value = '2026-02-14 06:20:02'
$ws.Range('J36').Value = '989.6 hPa'
$ws.Range('E37').Value = '2026-02-14 06:20:04'
$ws.Range('J37').Value = '991.4 hPa'
$ws.Range('N37').Value = '3.0 °C 5:57 TU'
$ws.Range('O37').Value = '3.4 °C'
$ws.Range('E38').Value = '2026-02-14 06:20:07'
$ws.Range('N38').Value = '6.4 °C 5:48 TU'
$ws.Range('O38').Value = '7.2 °C'
$ws.Range('E39').Value = '2026-02-14 06:20:09'
$ws.Range('I39').Value = '3.3 mm'
$ws.Range('N39').Value = '-5.8 °C 5:44 TU'
$ws.Range('O39').Value = '-5.0 °C'
$ws.Range('E40').Value = '2026-02-14 06:20:12'
$ws.Range('J40').Value = '993.0 hPa'
$ws.Range('M40').Value = '4.1 °C 5:58 TU'
$ws.Range('O40').Value = '2.5 °C'
$ws.Range('E41').Value = '2026-02-14 06:20:14'
$ws.Range('H41').NumberFormat = "@"
$ws.Range('H41').Value = '53%'
$ws.Range('J41').Value = '991.3 hPa'
$ws.Range('N41').Value = '10.8 °C 5:50 TU'
$ws.Range('O41').Value = '11.6 °C'
$ws.Range('E42').Value = '2026-02-14 06:20:17'
$ws.Range('O42').Value = '10.0 °C'
$ws.Range('E43').Value = '2026-02-14 06:20:19'
$ws.Range('N43').Value = '4.3 °C 5:55 TU'
$ws.Range('O43').Value = '5.8 °C'
$ws.Range('E44').Value = '2026-02-14 06:20:22'
$ws.Range('I44').Value = '14.3 mm'
$ws.Range('N44').Value = '-5.4 °C 5:58 TU'
$ws.Range('O44').Value = '-5.0 °C'
$ws.Range('E45').Value = '2026-02-14 06:20:24'
$ws.Range('I45').Value = '4.5 mm'
$ws.Range('J45').Value = '997.0 hPa'
$ws.Range('E46').Value = '2026-02-14 06:20:27'
$ws.Range('H46').NumberFormat = "@"
$ws.Range('H46').Value = '59%'
$ws.Range('J46').Value = '994.9 hPa'
$ws.Range('N46').Value = '9.2 °C 5:53 TU'
$ws.Range('O46').Value = '10.3 °C'
